# ---------------------------------------------------------------------------
# Add a new "2022-Q4" sheet (with fund holdings detail) right after the
# "总计" (summary) sheet, and update the summary sheet with the new
# 2022-Q4 totals row (shifting the existing quarters down by one row).
# ---------------------------------------------------------------------------

$xlPasteValues = -4163
$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by duplicating the "2022-Q3" sheet (so
#    that it starts out with the same column layout/formatting), placing it
#    immediately after "总计".
# ---------------------------------------------------------------------------
$q3Sheet.Copy([System.Reflection.Missing]::Value, $summarySheet)
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# The "2022-Q3" sheet only has 14 data rows (rows 2-15); "2022-Q4" needs 15
# data rows (rows 2-16), so extend the formatting of the last data row down
# to the new row 16 before filling in values.
$q4Sheet.Range("A15:H15").Copy()
$q4Sheet.Range("A16:H16").PasteSpecial($xlPasteFormats)

# Fund holdings detail for 2022-Q4.
# Columns: A=index, B=基金代码, C=基金名称, D=基金规模, E=股票总仓位,
#          F=仓位占比, G=持有市值(亿元), H=仓位排名
$q4Data = @(
    ,@(0,  '159869', '华夏中证动漫游戏ETF', '6.88', '99.30', '5.52', '0.3798', 6)
    ,@(1,  '320005', '诺安价值增长混合', '10.34', '83.51', '2.34', '0.2420', 10)
    ,@(2,  '516010', '国泰中证动漫游戏ETF', '4.38', '94.84', '5.29', '0.2317', 6)
    ,@(3,  '005585', '银河文体娱乐主题灵活配置混合A', '3.15', '88.90', '6.50', '0.2048', 3)
    ,@(4,  '001628', '招商体育文化休闲股票A', '2.33', '93.03', '4.65', '0.1083', 10)
    ,@(5,  '516770', '华泰柏瑞中证动漫游戏ETF', '1.10', '97.03', '5.13', '0.0564', 7)
    ,@(6,  '161036', '富国中证娱乐主题指数增强（LOF）A', '1.03', '93.82', '3.21', '0.0331', 8)
    ,@(7,  '015667', '银河文体娱乐主题灵活配置混合C', '0.38', '88.90', '6.50', '0.0247', 3)
    ,@(8,  '517500', '国泰中证沪港深动漫游戏ETF', '0.54', '98.42', '4.04', '0.0218', 7)
    ,@(9,  '004890', '中邮健康文娱灵活配置混合', '0.42', '92.60', '4.72', '0.0198', 4)
    ,@(10, '015395', '招商体育文化休闲股票C', '0.29', '93.03', '4.65', '0.0135', 10)
    ,@(11, '014246', '大摩现代服务业混合A', '0.17', '86.98', '5.85', '0.0099', 10)
    ,@(12, '014256', '富国中证娱乐主题指数增强（LOF）C', '0.13', '93.82', '3.21', '0.0042', 8)
    ,@(13, '005167', '嘉实润泽量化一年定期开放混合', '0.56', '27.25', '0.65', '0.0036', 2)
    ,@(14, '014247', '大摩现代服务业混合C', '0.06', '86.98', '5.85', '0.0035', 10)
)

$rowIndex = 2
foreach ($row in $q4Data) {
    $a = $row[0]; $b = $row[1]; $c = $row[2]; $d = $row[3]
    $e = $row[4]; $f = $row[5]; $g = $row[6]; $h = $row[7]

    # A and H are genuine numbers.
    $q4Sheet.Cells.Item($rowIndex, 1).Value = $a
    $q4Sheet.Cells.Item($rowIndex, 8).Value = $h

    # C (fund name) is plain text and never looks numeric, so it is safe to
    # assign directly.
    $q4Sheet.Cells.Item($rowIndex, 3).Value = $c

    # B, D, E, F, G all look numeric (fund codes / ratios) but must be
    # stored as *text*, matching the rest of the workbook. Writing them via
    # .Value would make Excel auto-convert them to numbers, so instead we
    # write a text formula and immediately convert it to a static value,
    # which keeps the original (non-numeric) cell formatting intact.
    $q4Sheet.Cells.Item($rowIndex, 2).Formula = '="' + $b + '"'
    $q4Sheet.Cells.Item($rowIndex, 4).Formula = '="' + $d + '"'
    $q4Sheet.Cells.Item($rowIndex, 5).Formula = '="' + $e + '"'
    $q4Sheet.Cells.Item($rowIndex, 6).Formula = '="' + $f + '"'
    $q4Sheet.Cells.Item($rowIndex, 7).Formula = '="' + $g + '"'

    $textRange = $q4Sheet.Range($q4Sheet.Cells.Item($rowIndex, 2), $q4Sheet.Cells.Item($rowIndex, 7))
    $textRange.Copy()
    $textRange.PasteSpecial($xlPasteValues)

    $rowIndex = $rowIndex + 1
}

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: push the existing quarters down one
#    row and insert the new 2022-Q4 totals at the top of the data.
# ---------------------------------------------------------------------------
# Row 6 does not exist yet, so copy column-A formatting from an existing
# data row before writing into it.
$summarySheet.Range("A2").Copy()
$summarySheet.Range("A6").PasteSpecial($xlPasteFormats)

$summarySheet.Cells.Item(6, 1).Value = 4
$summarySheet.Cells.Item(6, 2).Value = "2021-Q4"
$summarySheet.Cells.Item(6, 3).Value = 11
$summarySheet.Cells.Item(6, 4).Value = 0.86

$summarySheet.Cells.Item(5, 1).Value = 3
$summarySheet.Cells.Item(5, 2).Value = "2022-Q1"
$summarySheet.Cells.Item(5, 3).Value = 6
$summarySheet.Cells.Item(5, 4).Value = 0.59

$summarySheet.Cells.Item(4, 1).Value = 2
$summarySheet.Cells.Item(4, 2).Value = "2022-Q2"
$summarySheet.Cells.Item(4, 3).Value = 3
$summarySheet.Cells.Item(4, 4).Value = 0.48

$summarySheet.Cells.Item(3, 1).Value = 1
$summarySheet.Cells.Item(3, 2).Value = "2022-Q3"
$summarySheet.Cells.Item(3, 3).Value = 14
$summarySheet.Cells.Item(3, 4).Value = 1.04

$summarySheet.Cells.Item(2, 1).Value = 0
$summarySheet.Cells.Item(2, 2).Value = "2022-Q4"
$summarySheet.Cells.Item(2, 3).Value = 15
$summarySheet.Cells.Item(2, 4).Value = 1.36

# Restore the originally active sheet ("总计") as the active tab.
$summarySheet.Activate() | Out-Null
$summarySheet.Range("A1").Select() | Out-Null
